# Applies the numeric restatements + row-78 blank-out described in the
# commit "Adicionados balanc,os concatenados em uma unica planilha."
#
# Row 78 ("Part. de Acionistas Nao Controladores") had its P:AH cells
# (previously numeric 0s) turned into empty/blank cells once the sheet
# was concatenated with the other balance sheets. A handful of other
# cells across rows 58-61/63/67/69/73/74/79 got tiny (sub-cent) floating
# point restatements from the same concatenation/recalculation pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 58 ---
$ws.Range("H58").Value = 229406.944
$ws.Range("L58").Value = 358572.992
$ws.Range("P58").Value = 456922.016
$ws.Range("T58").Value = 46166.08
$ws.Range("X58").Value = 615957.12
$ws.Range("AF58").Value = 510714.048

# --- Row 59 ---
$ws.Range("H59").Value = -114373.008
$ws.Range("L59").Value = -156210.96
$ws.Range("T59").Value = 108250.016
$ws.Range("X59").Value = -418814.912
$ws.Range("AB59").Value = -472089.056

# --- Row 60 ---
$ws.Range("H60").Value = 115034.016
$ws.Range("L60").Value = 202362
$ws.Range("T60").Value = 154415.984
$ws.Range("AF60").Value = -213181.984

# --- Row 61 ---
$ws.Range("L61").Value = -86589
$ws.Range("P61").Value = -58017
$ws.Range("T61").Value = -96830.984
$ws.Range("X61").Value = -185326.992
$ws.Range("AB61").Value = -107957.992

# --- Row 63 ---
$ws.Range("AF63").Value = -53299

# --- Row 67 ---
$ws.Range("T67").Value = -28577.992
$ws.Range("AF67").Value = -68225

# --- Row 69 ---
$ws.Range("L69").Value = 115773
$ws.Range("P69").Value = 133854.008
$ws.Range("T69").Value = 57585.016
$ws.Range("AF69").Value = -363140.96

# --- Row 73 ---
$ws.Range("H73").Value = 4238
$ws.Range("P73").Value = 133984.024

# --- Row 74 ---
$ws.Range("H74").Value = -50303
$ws.Range("L74").Value = 17424.008
$ws.Range("P74").Value = 17603

# --- Row 78: P78:AH78 go from numeric 0 to blank/empty (text) cells ---
# A leading apostrophe forces Excel to store an explicit empty-text
# value (rather than fully clearing/removing the cell), which matches
# the target inline-string empty-cell state. That alone also marks the
# cells with a "quote prefix" style, so we then paste-format from the
# untouched neighbour O78 to restore the original (default) style
# while keeping the empty-text value.
$ws.Range("P78:AH78").Value = "'"
$ws.Range("O78").Copy()
$ws.Range("P78:AH78").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 79 ---
$ws.Range("H79").Value = 20282.008
$ws.Range("L79").Value = 85845
$ws.Range("P79").Value = 69474.008
$ws.Range("T79").Value = 30634.008
